$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:D1): rename columns to short english names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case fixes for state / municipality names ("de" -> "De", etc.) ---
$ws.Range("A26").Value  = "Ciudad De México"
$ws.Range("A34").Value  = "Estado De México"
$ws.Range("A38").Value  = "Guanajuato"
$ws.Range("B38").Value  = "Apaseo El Alto"
$ws.Range("B40").Value  = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B42").Value  = "Purísima Del Rincón"
$ws.Range("B46").Value  = "Acapulco De Juárez"
$ws.Range("B50").Value  = "Iguala De La Independencia"
$ws.Range("B51").Value  = "Zihuatanejo De Azueta"
$ws.Range("B52").Value  = "Técpan De Galeana"
$ws.Range("B63").Value  = "Tamazula De Gordiano"
$ws.Range("B64").Value  = "Tlajomulco De Zúñiga"
$ws.Range("B65").Value  = "Unión De Tula"
$ws.Range("B69").Value  = "Coalcomán De Vázquez Pallares"
$ws.Range("B84").Value  = "Mier Y Noriega"
$ws.Range("B85").Value  = "San Nicolás De Los Garza"
$ws.Range("B91").Value  = "Zapotitlán Del Río"
$ws.Range("B101").Value = "Ciudad Del Maíz"
$ws.Range("B104").Value = "Villa De Ramos"
$ws.Range("B123").Value = "Martínez De La Torre"
$ws.Range("B129").Value = "Jiménez Del Teul"
$ws.Range("B131").Value = "Villa De Cos"

# --- Minor float re-computation of a percentage cell ---
$ws.Range("D54").Value = 0.09803921568627452

# --- Drop the trailing footer/metadata rows (135-139) ---
$ws.Rows("135:139").Delete()

Write-Host "edit complete"
